# Apply crypto price/volume updates per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = "61.738.14"
$ws.Range("E2").Value2 = "  +2.34%  "
$ws.Range("D3").Value2 = "2.382.63"
$ws.Range("E3").Value2 = "  +1.48%  "
$ws.Range("E4").Value2 = "  -0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value2 = "553.50"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value2 = "  +2.41%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value2 = "141.22"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value2 = "  +3.42%  "
$ws.Range("E7").Value2 = "  -0.16%  "
$ws.Range("E8").Value2 = "  +0.66%  "
$ws.Range("D9").Value2 = "2.383.75"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value2 = "0.109"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value2 = "  +5.03%  "
$ws.Range("E11").Value2 = "  +2.19%  "
$ws.Range("E12").Value2 = "  +2.36%  "
$ws.Range("E13").Value2 = "  +4.11%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value2 = "25.71"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value2 = "  +5.86%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value2 = "0.0000176"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value2 = "  +10.19%  "
$ws.Range("D16").Value2 = "2.810.67"
$ws.Range("E16").Value2 = "  +1.29%  "
$ws.Range("D17").Value2 = "61.598.59"
$ws.Range("E17").Value2 = "  +2.01%  "
$ws.Range("D18").Value2 = "2.384.49"
$ws.Range("E18").Value2 = "  +1.36%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value2 = "11.02"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value2 = "  +4.81%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value2 = "323.75"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value2 = "  +4.10%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value2 = "4.18"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value2 = "  +3.17%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value2 = "6.71"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value2 = "  +2.24%  "
$ws.Range("E23").Value2 = "  +0.00%  "
$ws.Range("E24").Value2 = "  -4.26%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value2 = "64.37"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value2 = "  +2.31%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value2 = "8.93"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value2 = "  +4.72%  "
$ws.Range("E27").Value2 = "  -0.26%  "
$ws.Range("B28").Value2 = "Bittensor"
$ws.Range("C28").Value2 = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value2 = "536.58"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value2 = "  +7.84%  "
$ws.Range("B29").Value2 = "WrappedeETH"
$ws.Range("C29").Value2 = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D29").Value2 = "2.498.04"
$ws.Range("E29").Value2 = "  +1.16%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value2 = "8.29"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value2 = "  +4.89%  "
$ws.Range("D31").Value2 = "0.0₃0923"
$ws.Range("E31").Value2 = "  +4.99%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value2 = "1.42"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value2 = "  +3.30%  "
$ws.Range("E33").Value2 = "  +3.75%  "
$ws.Range("E34").Value2 = "  +4.39%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value2 = "1.53"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value2 = "  +1.55%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value2 = "5.74"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value2 = "  +10.50%  "
$ws.Range("E37").Value2 = "  -0.13%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value2 = "4.76"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value2 = "  +4.99%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value2 = "1.93"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value2 = "  +9.26%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value2 = "0.382"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value2 = "  +3.41%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value2 = "18.59"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value2 = "  +1.76%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value2 = "146.89"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value2 = "  +7.69%  "
$ws.Range("E43").Value2 = "  -0.03%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value2 = "41.53"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value2 = "  +3.93%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value2 = "148.91"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value2 = "  +5.68%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value2 = "2.21"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value2 = "  +5.65%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value2 = "3.62"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value2 = "  +2.98%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value2 = "0.0529"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value2 = "  +4.49%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value2 = "20.17"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value2 = "  +4.63%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value2 = "0.586"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value2 = "  +3.48%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value2 = "0.0907"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value2 = "  +1.63%  "
